# Homework-aAlphaBio.pptx edit:
#  1) Move the "Relevant literature, code, and data" (BERT) slide from position 2
#     to position 6 (i.e. after the "aa sequence" slide, right before "Results").
#  2) Re-purpose/re-format that slide as a "BERT Stuff" recap slide:
#       - widen + retitle the heading textbox
#       - reflow + shrink the bibliography textbox, 16pt body text
#       - add a new note textbox under the heading
#
$p = $ppt.ActivePresentation

# --- Step 1: locate & move the BERT slide -------------------------------
$bertSlide = $p.Slides.Item(2)
$bertSlide.MoveTo(6)

# Re-fetch by position now that the deck has been reordered.
$s = $p.Slides.Item(6)

# --- Step 2a: heading textbox (TextBox 1) -------------------------------
$title = $s.Shapes.Item(1)
$title.TextFrame.TextRange.Text = "BERT Stuff:   Relevant literature, code, and data"
$title.Width = 6376297 / 12700

# --- Step 2b: bibliography textbox (TextBox 2) --------------------------
$body = $s.Shapes.Item(2)
$body.Left = 765387 / 12700
$body.Top = 1291616 / 12700
$body.Width = 9574610 / 12700
$body.Height = 5016758 / 12700

$bodyRange = $body.TextFrame.TextRange
$paraCount = $bodyRange.Paragraphs().Count
for ($i = 1; $i -le $paraCount; $i++) {
    $para = $bodyRange.Paragraphs($i, 1)
    $para.Font.Size = 16
}

# --- Step 2c: new note textbox (TextBox 4) ------------------------------
# Burn an id so the new shape lands on id=5 / name "TextBox 4", matching the
# id the real edit session left behind (a shape was created & removed earlier
# in that session).
$scratch = $s.Shapes.AddTextbox(1, 0, 0, 10, 10)
$scratch.Delete()

$note = $s.Shapes.AddTextbox(1, 140127 / 12700, 679781 / 12700, 9369296 / 12700, 307777 / 12700)
$note.TextFrame.WordWrap = 0
$note.TextFrame.AutoSize = 1
$note.Fill.Visible = 0
$note.TextFrame.TextRange.Text = "Note: I did most of the BERT coding/testing in the 2 weeks before my call with Adrian (i.e. before getting the HW problem)"
$note.TextFrame.TextRange.Font.Size = 14

Write-Output "done"
